# Insert two new weekly price rows for "Brócoli" above the existing
# row 564 (pushing the prior rows 564-612 down to 566-614), then
# populate the two new rows with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 564:612 -> 566:614 by inserting two blank rows at 564.
$ws.Rows("564:565").Insert()

# New row 564: Primera quality entry for 2023-09-25 (serial 45194).
$ws.Range("A564").Value = 7
$ws.Range("B564").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C564").Value = "Ñuble"
$ws.Range("D564").Value = 45194
$ws.Range("E564").Value = 16
$ws.Range("F564").Value = 100112023
$ws.Range("G564").Value = "Brócoli"
$ws.Range("H564").Value = "Sin especificar"
$ws.Range("I564").Value = "Primera"
$ws.Range("J564").Value = 300
$ws.Range("K564").Value = 1000
$ws.Range("L564").Value = 1000
$ws.Range("M564").Value = 1000
$ws.Range("N564").Value = "$/unidad"
$ws.Range("O564").Value = "Región del Maule"
$ws.Range("P564").Value = 1000
$ws.Range("Q564").Value = 1
$ws.Range("R564").Value = "Hortaliza"

# New row 565: Segunda quality entry for 2023-09-25 (serial 45194).
$ws.Range("A565").Value = 7
$ws.Range("B565").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C565").Value = "Ñuble"
$ws.Range("D565").Value = 45194
$ws.Range("E565").Value = 16
$ws.Range("F565").Value = 100112023
$ws.Range("G565").Value = "Brócoli"
$ws.Range("H565").Value = "Sin especificar"
$ws.Range("I565").Value = "Segunda"
$ws.Range("J565").Value = 400
$ws.Range("K565").Value = 800
$ws.Range("L565").Value = 800
$ws.Range("M565").Value = 800
$ws.Range("N565").Value = "$/unidad"
$ws.Range("O565").Value = "Región del Maule"
$ws.Range("P565").Value = 800
$ws.Range("Q565").Value = 1
$ws.Range("R565").Value = "Hortaliza"
